$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.941.75'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '1.643.79'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''217.43'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = '''0.5234'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').Value = '''1.003'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = '''0.2614'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('D9').Value = '''0.06269'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = '''20.37'
$ws.Range('E10').Value = '  -3.61%  '
$ws.Range('D11').Value = '''0.07736'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '''4.444'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.651.66'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').Value = '''0.5431'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').Value = '0.0₅8053'
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').Value = '''64.66'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = '25.968.24'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').Value = '''4.543'
$ws.Range('E19').Value = '  -2.87%  '
$ws.Range('D20').Value = '''191.57'
$ws.Range('D21').Value = '''10.02'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = '''5.964'
$ws.Range('E23').Value = '  -2.43%  '
$ws.Range('D24').Value = '''139.92'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('D25').Value = '''0.1235'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').Value = '''7.244'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '''16.15'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '''1.420'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').Value = '''0.05921'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('D31').Value = '''3.484'
$ws.Range('E31').Value = '  -1.97%  '
$ws.Range('D32').Value = '''3.230'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').Value = '''1.524'
$ws.Range('E33').Value = '  -7.88%  '
$ws.Range('D34').Value = '''2.415'
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').Value = '''0.9389'
$ws.Range('E35').Value = '  -4.39%  '
$ws.Range('D36').Value = '''2.740'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('D37').Value = '''0.5724'
$ws.Range('E37').Value = '  -3.44%  '
$ws.Range('D38').Value = '''0.01604'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').Value = '''5.848'
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('D40').Value = '''0.8460'
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = '''100.45'
$ws.Range('E42').Value = '  +0.53%  '
$ws.Range('D43').Value = '1.000.07'
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('D44').Value = '1.785.59'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').Value = '''56.44'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('D47').Value = '''1.005'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').Value = '''0.4288'
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('D49').Value = '''1.473'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').Value = '''0.05148'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').Value = '''7.836'
$ws.Range('E51').Value = '  -3.73%  '
